$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-12 from
# 45183 (2023-09-14) to 45184 (2023-09-15).
$ws.Range("C2:C12").Value = 45184
